$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make sure the date column's number format is set before any value is
# written into it, so the workbook doesn't mint a spurious custom numFmt
# entry (Excel recognizes "m/d/yy h:mm" as the built-in numFmtId 22).
$ws.Range("A1:A5").NumberFormat = "m/d/yy h:mm"

$rows = @(
    @(42606.574745370373, 28),
    @(42606.575844907406, 42),
    @(42606.580891203703, 31),
    @(42606.581284722219, 53)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $rows[$i][0]
    $ws.Cells.Item($r, 2).Value = $rows[$i][1]
    for ($c = 3; $c -le 13; $c++) {
        $ws.Cells.Item($r, $c).Value = 0
    }
    $ws.Cells.Item($r, 14).Value = "Random"
}

# Widen column A to fit the new date/time values (closest width reachable
# through the ColumnWidth property's internal rounding).
$ws.Columns(1).ColumnWidth = 14
